$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Statistics" (first sheet): refresh rows 2-14 with new
# simulation readings and drop the now-unused trailing rows 15-19.
# ---------------------------------------------------------------
$stats = $wb.Worksheets.Item("Statistics")

$statsData = @(
    @("2024-08-28 17:21:18", 37.55490894322983, 8),
    @("2024-08-28 17:21:20", 36.15062400517214, 12),
    @("2024-08-28 17:21:22", 35.15737758284577, 17),
    @("2024-08-28 17:21:24", 33.27148322901341, 20),
    @("2024-08-28 17:21:26", 35.68417963447074, 26),
    @("2024-08-28 17:21:29", 34.05210373195279, 27),
    @("2024-08-28 17:21:31", 35.11168667382415, 29),
    @("2024-08-28 17:21:33", 34.45730914588974, 30),
    @("2024-08-28 17:21:35", 30.08983503040512, 34),
    @("2024-08-28 17:21:37", 28.57848323067096, 35),
    @("2024-08-28 17:21:39", 31.11678594657995, 36),
    @("2024-08-28 17:21:41", 28.00793533950663, 38),
    @("2024-08-28 17:21:43", 20.17901844736679, 38)
)

for ($i = 0; $i -lt $statsData.Count; $i++) {
    $r = $i + 2
    $row = $statsData[$i]
    $stats.Cells.Item($r, 1).Value = $row[0]
    $stats.Cells.Item($r, 2).Value = $row[1]
    $stats.Cells.Item($r, 3).Value = $row[2]
}

# Remove the rows that no longer exist in the refreshed data set
# (work from the bottom up so row numbers stay stable).
for ($r = 19; $r -ge 15; $r--) {
    $stats.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------
# Sheet "Accidents" (second sheet): refresh rows 2-5 with new
# accident records and drop the now-unused trailing rows 6-14.
# ---------------------------------------------------------------
$accidents = $wb.Worksheets.Item("Accidents")

$accidentsData = @(
    @("2024-08-28 17:21:33", "Truck and Car", "31.04 and 38.58", 1),
    @("2024-08-28 17:21:41", "Car and Truck", "35.93 and 37.05", 2),
    @("2024-08-28 17:21:42", "Car and Car", "0.00 and 33.31", 2),
    @("2024-08-28 17:21:44", "Car and Car", "0.00 and 37.54", 2)
)

for ($i = 0; $i -lt $accidentsData.Count; $i++) {
    $r = $i + 2
    $row = $accidentsData[$i]
    $accidents.Cells.Item($r, 1).Value = $row[0]
    $accidents.Cells.Item($r, 2).Value = $row[1]
    $accidents.Cells.Item($r, 3).Value = $row[2]
    $accidents.Cells.Item($r, 4).Value = $row[3]
}

for ($r = 14; $r -ge 6; $r--) {
    $accidents.Rows.Item($r).Delete()
}
